$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A39").Value = "GRT-USD"
